# Updates the cryptos list sheet with latest prices / volume changes.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellRef, $text) {
    $rng = $ws.Range($cellRef)
    $rng.NumberFormat = "@"
    $rng.Value = $text
    $rng.Style = "Normal"
}

# Row 2 - Bitcoin
Set-TextValue "D2" "66.820.72"
$ws.Range("E2").Value = "  -0.60%  "

# Row 3 - Ethereum
Set-TextValue "D3" "3.452.32"
$ws.Range("E3").Value = "  -1.51%  "

# Row 4 - TetherUSD
$ws.Range("E4").Value = "  +0.03%  "

# Row 5 - BNB
Set-TextValue "D5" "591.53"
$ws.Range("E5").Value = "  -1.49%  "

# Row 6 - Solana
Set-TextValue "D6" "178.21"
$ws.Range("E6").Value = "  +2.17%  "

# Row 7 - XRP
Set-TextValue "D7" "0.609"
$ws.Range("E7").Value = "  +3.93%  "

# Row 8 - USDC
$ws.Range("E8").Value = "  +0.03%  "

# Row 9 - LidoStakedEther
Set-TextValue "D9" "3.451.05"
$ws.Range("E9").Value = "  -1.51%  "

# Row 10 - Dogecoin
Set-TextValue "D10" "0.139"
$ws.Range("E10").Value = "  +5.57%  "

# Row 12 - Cardano
$ws.Range("E12").Value = "  -0.13%  "

# Row 13 - WrappedliquidstakedEther2.0
Set-TextValue "D13" "4.047.25"
$ws.Range("E13").Value = "  -1.57%  "

# Row 14 - Avalanche
Set-TextValue "D14" "31.66"
$ws.Range("E14").Value = "  +3.40%  "

# Row 15 - TRON
$ws.Range("E15").Value = "  -0.49%  "

# Row 16 - WrappedBTC
Set-TextValue "D16" "66.858.98"
$ws.Range("E16").Value = "  -0.49%  "

# Row 17 - ShibaInu
Set-TextValue "D17" "0.0000177"
$ws.Range("E17").Value = "  -1.26%  "

# Row 18 - WrappedEther
Set-TextValue "D18" "3.451.59"
$ws.Range("E18").Value = "  -1.73%  "

# Row 19 - Polkadot
$ws.Range("E19").Value = "  -1.29%  "

# Row 20 - Chainlink
Set-TextValue "D20" "14.13"
$ws.Range("E20").Value = "  -3.31%  "

# Row 21 - BitcoinCash
Set-TextValue "D21" "389.01"
$ws.Range("E21").Value = "  -1.24%  "

# Row 22 - Uniswap
$ws.Range("E22").Value = "  -1.22%  "

# Row 23 - Dai
$ws.Range("E23").Value = "  -0.20%  "

# Row 25 - Litecoin
Set-TextValue "D25" "71.88"
$ws.Range("E25").Value = "  -2.01%  "

# Row 26 - Polygon
Set-TextValue "D26" "0.534"
$ws.Range("E26").Value = "  -0.63%  "

# Row 27 - PEPE
$ws.Range("E27").Value = "  -0.72%  "

# Row 28 - InternetComputer(DFINITY)
Set-TextValue "D28" "10.23"
$ws.Range("E28").Value = "  +0.77%  "

# Row 29 - Kaspa
Set-TextValue "D29" "0.174"
$ws.Range("E29").Value = "  -3.70%  "

# Row 30 - Binance-PegBSC-USD
Set-TextValue "D30" "1.00"
$ws.Range("E30").Value = "  +0.56%  "

# Row 31 - NEARProtocol
Set-TextValue "D31" "6.14"
$ws.Range("E31").Value = "  -0.43%  "

# Row 32 - now Fetch.AI (was PancakeSwap)
$ws.Range("B32").Value = "Fetch.AI"
$ws.Range("C32").Value = "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
Set-TextValue "D32" "1.40"
$ws.Range("E32").Value = "  -1.91%  "

# Row 33 - now PancakeSwap (was Fetch.AI)
$ws.Range("B33").Value = "PancakeSwap"
$ws.Range("C33").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
Set-TextValue "D33" "2.05"
$ws.Range("E33").Value = "  -0.69%  "

# Row 34 - EthereumClassic
Set-TextValue "D34" "23.38"
$ws.Range("E34").Value = "  -1.15%  "

# Row 35 - Aptos
Set-TextValue "D35" "7.29"

# Row 36 - USDe
$ws.Range("E36").Value = "  -0.07%  "

# Row 37 - ImmutableX
$ws.Range("E37").Value = "  -3.38%  "

# Row 38 - Monero
Set-TextValue "D38" "163.53"
$ws.Range("E38").Value = "  -0.26%  "

# Row 39 - Mantle
Set-TextValue "D39" "0.874"
$ws.Range("E39").Value = "  -0.52%  "

# Row 40 - dogwifhat
Set-TextValue "D40" "2.79"
$ws.Range("E40").Value = "  +9.85%  "

# Row 41 - Stacks
$ws.Range("E41").Value = "  -2.75%  "

# Row 42 - RenderToken
Set-TextValue "D42" "6.75"
$ws.Range("E42").Value = "  -3.88%  "

# Row 43 - Filecoin
Set-TextValue "D43" "4.64"
$ws.Range("E43").Value = "  -0.85%  "

# Row 44 - EnergySwap
Set-TextValue "D44" "26.10"
$ws.Range("E44").Value = "  -0.01%  "

# Row 45 - Hedera
Set-TextValue "D45" "0.0716"
$ws.Range("E45").Value = "  -2.20%  "

# Row 46 - Maker
Set-TextValue "D46" "2.730.58"
$ws.Range("E46").Value = "  -2.69%  "

# Row 47 - InjectiveProtocol
Set-TextValue "D47" "26.02"
$ws.Range("E47").Value = "  -5.76%  "

# Row 48 - OKB
$ws.Range("E48").Value = "  -3.35%  "

# Row 49 - VeChain
$ws.Range("E49").Value = "  -1.74%  "

# Row 50 - Bittensor
Set-TextValue "D50" "324.83"
$ws.Range("E50").Value = "  -4.80%  "

# Row 51 - ONDO
$ws.Range("E51").Value = "  -3.85%  "
